$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1.83
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 2.08
$ws.Range("H3").Value = 3.9
$ws.Range("J3").Value = 3.25
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.39
$ws.Range("P3").Value = 1.68
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.25
$ws.Range("T3").Value = 1.94
$ws.Range("U3").Value = 1.84
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 1.92
$ws.Range("AF3").Value = 13.5
$ws.Range("AG3").Value = 13
$ws.Range("AJ3").Value = 28
$ws.Range("AK3").Value = 29
$ws.Range("AL3").Value = 55
$ws.Range("AN3").Value = 22
$ws.Range("F4").Value = 1.85
$ws.Range("H4").Value = 4.3
$ws.Range("N4").Value = 3.65
$ws.Range("P4").Value = 1.92
$ws.Range("Q4").Value = 1.94
$ws.Range("Z4").Value = 36
$ws.Range("AC4").Value = 9.800000000000001
$ws.Range("G5").Value = 2.66
$ws.Range("J5").Value = 3.25
$ws.Range("N5").Value = 3.4
$ws.Range("T5").Value = 1.65
$ws.Range("U5").Value = 1.93
$ws.Range("V5").Value = 1.37
$ws.Range("W5").Value = 1.62
$ws.Range("R6").Value = 1.69
$ws.Range("S6").Value = 2.18
$ws.Range("T6").Value = 1.65
$ws.Range("U6").Value = 1.9
$ws.Range("J7").Value = 6.8
$ws.Range("P7").Value = 2.64
$ws.Range("U7").Value = 1.83
$ws.Range("G8").Value = 2.52
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 3.3
$ws.Range("Q8").Value = 2.24
$ws.Range("U8").Value = 2.02
$ws.Range("W8").Value = 1.66
$ws.Range("AN8").Value = 26
$ws.Range("G9").Value = 1.59
$ws.Range("H9").Value = 6.2
$ws.Range("I9").Value = 7.6
$ws.Range("K9").Value = 5.6
$ws.Range("N9").Value = 5.1
$ws.Range("S9").Value = 2.32
$ws.Range("T9").Value = 1.71
$ws.Range("U9").Value = 2.18
$ws.Range("W9").Value = 2.68
$ws.Range("X9").Value = 28
$ws.Range("AH9").Value = 23
$ws.Range("AN9").Value = 6
$ws.Range("S10").Value = 2.4
$ws.Range("U10").Value = 1.79
$ws.Range("Z10").Value = 85
$ws.Range("H11").Value = 1.68
$ws.Range("K11").Value = 4.3
$ws.Range("O11").Value = 1.34
$ws.Range("P11").Value = 1.83
$ws.Range("Q11").Value = 1.98
$ws.Range("T11").Value = 1.94
$ws.Range("U11").Value = 1.89
$ws.Range("X11").Value = 14.5
$ws.Range("Y11").Value = 8
$ws.Range("AC11").Value = 9.4
$ws.Range("AF11").Value = 50
$ws.Range("AG11").Value = 23
$ws.Range("AH11").Value = 25
$ws.Range("AI11").Value = 44
$ws.Range("AO11").Value = 12.5
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = 9.800000000000001
$ws.Range("H12").Value = 1.37
$ws.Range("K12").Value = 6.2
$ws.Range("L12").Value = 1.23
$ws.Range("N12").Value = 5.6
$ws.Range("P12").Value = 2.52
$ws.Range("Q12").Value = 1.48
$ws.Range("R12").Value = 1.61
$ws.Range("S12").Value = 2.28
$ws.Range("T12").Value = 1.81
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 3.35
$ws.Range("X12").Value = 34
$ws.Range("Z12").Value = 970
$ws.Range("AA12").Value = 15
$ws.Range("AB12").Value = 36
$ws.Range("AG12").Value = 38
$ws.Range("AH12").Value = 28
$ws.Range("AI12").Value = 34
$ws.Range("AO12").Value = 5.3
$ws.Range("U13").Value = 1.62
$ws.Range("AC13").Value = 9.6
$ws.Range("L14").Value = 1.16
$ws.Range("P14").Value = 1.56
$ws.Range("R14").Value = 1.54
$ws.Range("Q15").Value = 2.24
$ws.Range("AO15").Value = 32
$ws.Range("F16").Value = 1.58
$ws.Range("G16").Value = 1.75
$ws.Range("J16").Value = 3.4
$ws.Range("Q16").Value = 1.89
$ws.Range("I17").Value = 2.32
$ws.Range("S17").Value = 4.8
$ws.Range("T17").Value = 2.04
$ws.Range("U17").Value = 1.79
$ws.Range("V17").Value = 1.79
$ws.Range("X18").Value = 15
$ws.Range("AB18").Value = 6.6
$ws.Range("AN18").Value = 12.5
$ws.Range("G19").Value = 1.51
$ws.Range("I19").Value = 8.800000000000001
$ws.Range("Q19").Value = 1.73
$ws.Range("W19").Value = 2.96
$ws.Range("H20").Value = 2.02
$ws.Range("I20").Value = 2.14
$ws.Range("J20").Value = 3.35
$ws.Range("P20").Value = 1.68
$ws.Range("T20").Value = 1.8
$ws.Range("V20").Value = 1.87
$ws.Range("Y20").Value = 8.199999999999999
$ws.Range("Z20").Value = 13
$ws.Range("AA20").Value = 28
$ws.Range("AB20").Value = 15.5
$ws.Range("AC20").Value = 9
$ws.Range("AD20").Value = 13
$ws.Range("AF20").Value = 34
$ws.Range("AG20").Value = 21
$ws.Range("AK20").Value = 80
$ws.Range("AL20").Value = 95
$ws.Range("AO20").Value = 23
$ws.Range("J21").Value = 3.25
$ws.Range("O21").Value = 1.44
$ws.Range("P21").Value = 1.71
$ws.Range("Q21").Value = 2.38
$ws.Range("AO21").Value = 70
$ws.Range("F22").Value = 1.91
$ws.Range("G22").Value = 1.98
$ws.Range("H22").Value = 4.1
$ws.Range("P22").Value = 2.04
$ws.Range("Q22").Value = 1.78
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("W22").Value = 2.02
$ws.Range("AF22").Value = 12.5
$ws.Range("AG22").Value = 10
$ws.Range("F23").Value = 1.87
$ws.Range("G23").Value = 1.88
$ws.Range("S23").Value = 2.5
$ws.Range("W23").Value = 2.12
$ws.Range("F24").Value = 2.72
$ws.Range("G24").Value = 2.76
$ws.Range("H24").Value = 2.82
$ws.Range("I24").Value = 2.86
$ws.Range("J24").Value = 3.45
$ws.Range("O24").Value = 1.35
$ws.Range("W24").Value = 1.56
$ws.Range("AF24").Value = 17.5
$ws.Range("AI24").Value = 44
$ws.Range("AJ24").Value = 42
$ws.Range("AN24").Value = 26
$ws.Range("F25").Value = 1.19
$ws.Range("G25").Value = 1.2
$ws.Range("H25").Value = 19
$ws.Range("I25").Value = 23
$ws.Range("J25").Value = 8.4
$ws.Range("K25").Value = 9
$ws.Range("N25").Value = 5.2
$ws.Range("R25").Value = 1.57
$ws.Range("V25").Value = 1.04
$ws.Range("W25").Value = 6
$ws.Range("AD25").Value = 90
$ws.Range("AH25").Value = 55
$ws.Range("F26").Value = 3.15
$ws.Range("G26").Value = 3.65
$ws.Range("H26").Value = 2.06
$ws.Range("I26").Value = 2.24
$ws.Range("J26").Value = 4.1
$ws.Range("K26").Value = 4.6
$ws.Range("L26").Value = 1.23
$ws.Range("N26").Value = 5.2
$ws.Range("O26").Value = 1.19
$ws.Range("P26").Value = 2.44
$ws.Range("Q26").Value = 1.57
$ws.Range("R26").Value = 1.6
$ws.Range("S26").Value = 2.36
$ws.Range("T26").Value = 1.58
$ws.Range("U26").Value = 2.48
$ws.Range("V26").Value = 1.8
$ws.Range("W26").Value = 1.38
$ws.Range("X26").Value = 34
$ws.Range("Z26").Value = 21
$ws.Range("AA26").Value = 32
$ws.Range("AC26").Value = 10.5
$ws.Range("AE26").Value = 21
$ws.Range("AG26").Value = 15
$ws.Range("AJ26").Value = 70
$ws.Range("AK26").Value = 36
$ws.Range("AO26").Value = 11
